# The deck ships two DrawingML theme parts: ppt/theme/theme1.xml (plain
# default "Office Theme" palette, wired only to the Notes Master) and
# ppt/theme/theme2.xml (the "Integral"/"Red Violet" palette, wired to the
# Slide Master and to the presentation's own theme relationship, i.e. the
# theme that is actually "live"/visible on every slide). The commit swaps
# the contents of the two parts so the live theme becomes the plain Office
# palette. Re-colour the live theme (reached through any slide's
# ThemeColorScheme, since there is a single Slide Master/Design in this
# deck) to the Office defaults, matching the colours that theme1.xml
# already carries.

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Colors(1).RGB  = RGBVal 0   0   0     # dk1      000000
$tcs.Colors(2).RGB  = RGBVal 255 255 255   # lt1      FFFFFF
$tcs.Colors(3).RGB  = RGBVal 68  84  106   # dk2      44546A
$tcs.Colors(4).RGB  = RGBVal 231 230 230   # lt2      E7E6E6
$tcs.Colors(5).RGB  = RGBVal 91  155 213   # accent1  5B9BD5
$tcs.Colors(6).RGB  = RGBVal 237 125 49    # accent2  ED7D31
$tcs.Colors(7).RGB  = RGBVal 165 165 165   # accent3  A5A5A5
$tcs.Colors(8).RGB  = RGBVal 255 192 0     # accent4  FFC000
$tcs.Colors(9).RGB  = RGBVal 68  114 196   # accent5  4472C4
$tcs.Colors(10).RGB = RGBVal 112 173 71    # accent6  70AD47
$tcs.Colors(11).RGB = RGBVal 5   99  193   # hlink    0563C1
$tcs.Colors(12).RGB = RGBVal 149 79  114   # folHlink 954F72
